$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2024-02-27 Tuesday" "2024-02-28 Wednesday"

Replace-Text "967×2=1934" "838×6=5028"
Replace-Text "516×2=1032" "175×4=700"
Replace-Text "913×5=4565" "531×2=1062"
Replace-Text "632×7=4424" "153×5=765"
Replace-Text "764×6=4584" "728×7=5096"
Replace-Text "746×7=5222" "521×9=4689"
Replace-Text "479×2=958" "804×6=4824"
Replace-Text "780×8=6240" "557×4=2228"
Replace-Text "209×4=836" "965×7=6755"
Replace-Text "475×8=3800" "470×9=4230"
Replace-Text "647×9=5823" "365×9=3285"
Replace-Text "102×7=714" "150×7=1050"
Replace-Text "310×3=930" "559×2=1118"
Replace-Text "801×7=5607" "876×6=5256"
Replace-Text "863×2=1726" "340×2=680"
Replace-Text "941×4=3764" "376×8=3008"
Replace-Text "640×6=3840" "846×5=4230"
Replace-Text "396×9=3564" "861×3=2583"
Replace-Text "898×7=6286" "227×9=2043"
Replace-Text "734×8=5872" "358×5=1790"
Replace-Text "862×4=3448" "559×6=3354"
Replace-Text "288×2=576" "342×3=1026"
Replace-Text "567×5=2835" "855×7=5985"
Replace-Text "143×9=1287" "233×4=932"
Replace-Text "841×9=7569" "159×9=1431"
